# Apply updated cryptocurrency price/volume figures (and the Stellar /
# InjectiveProtocol row swap) to match the refreshed "cryptos" feed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new value. Cells whose new text would otherwise be auto-parsed
# by Excel as a plain number (e.g. "601.53") are written with a leading
# apostrophe to force a Text entry, matching the workbook's existing
# inline-string ("Price" column looks numeric but is stored as text)
# convention; the quote-prefix style that creates is then reset back to
# "Normal" so the cell keeps its original (unstyled) appearance.
$updates = [ordered]@{
    'D2' = '63.960.06'
    'E2' = '  -1.71%  '
    'D3' = '3.144.60'
    'E3' = '  -1.36%  '
    'E4' = '  -0.01%  '
    'D5' = '601.53'
    'E5' = '  -2.35%  '
    'D6' = '141.98'
    'E6' = '  -4.11%  '
    'E7' = '  +0.07%  '
    'D8' = '3.139.35'
    'E8' = '  -1.46%  '
    'D9' = '0.527'
    'E9' = '  -0.83%  '
    'E10' = '  -3.34%  '
    'E11' = '  -2.57%  '
    'E12' = '  -2.71%  '
    'E13' = '  -4.37%  '
    'D14' = '34.85'
    'E14' = '  -3.68%  '
    'D15' = '3.664.29'
    'E15' = '  -1.34%  '
    'D17' = '63.926.05'
    'E17' = '  -1.83%  '
    'D18' = '3.137.38'
    'E18' = '  -1.44%  '
    'E19' = '  -2.12%  '
    'D20' = '486.49'
    'E20' = '  +0.27%  '
    'D21' = '14.67'
    'E21' = '  -0.99%  '
    'E22' = '  -2.29%  '
    'D23' = '7.73'
    'E23' = '  -3.36%  '
    'D24' = '88.43'
    'E24' = '  +4.30%  '
    'E25' = '  -5.56%  '
    'E26' = '  -0.01%  '
    'E27' = '  -2.93%  '
    'E28' = '  -6.81%  '
    'E29' = '  -2.59%  '
    'E30' = '  -3.72%  '
    'D31' = '27.41'
    'E31' = '  +1.93%  '
    'E32' = '  -8.18%  '
    'E33' = '  +0.00%  '
    'E34' = '  -3.73%  '
    'E35' = '  -3.27%  '
    'E36' = '  -0.60%  '
    'D37' = '52.71'
    'E37' = '  -1.19%  '
    'D38' = '0.0₃0742'
    'E38' = '  -7.41%  '
    'D39' = '2.91'
    'E39' = '  -9.90%  '
    'E40' = '  -1.67%  '
    'D41' = '430.91'
    'E41' = '  -8.58%  '
    'E42' = '  -0.85%  '
    'E43' = '  -0.95%  '
    'D44' = '2.910.74'
    'E44' = '  +1.22%  '
    'E45' = '  -4.91%  '
    'D46' = '2.18'
    'E46' = '  -7.53%  '
    'E47' = '  -3.37%  '
    'E48' = '  -0.06%  '
    'B49' = 'InjectiveProtocol'
    'C49' = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
    'D49' = '25.71'
    'E49' = '  -4.97%  '
    'B50' = 'Stellar'
    'C50' = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
    'D50' = '0.114'
    'E50' = '  -0.45%  '
    'D51' = '121.03'
    'E51' = '  -0.19%  '
}

$forceText = @(
    'D5'
    'D6'
    'D9'
    'D14'
    'D20'
    'D21'
    'D23'
    'D24'
    'D31'
    'D37'
    'D39'
    'D41'
    'D46'
    'D49'
    'D50'
    'D51'
)

foreach ($cell in $updates.Keys) {
    $value = $updates[$cell]
    $range = $ws.Range($cell)
    if ($forceText -contains $cell) {
        $range.Value = "'" + $value
        $range.Style = 'Normal'
    } else {
        $range.Value = $value
    }
}
